$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New GUID Directory rows (AllHabits heat-map captures), appended after the
# existing 61 rows (header + 60 data rows -> new data rows 62-67).
# Column A values are zero-padded numeric-looking codes ("000061" ...)
# that must be stored as literal text (leading zeros preserved), matching
# the existing rows in the sheet. Forcing a Text number format for the
# assignment keeps the input parser from collapsing them to the integer
# 61/62/etc, then resetting the style back to Normal/General afterwards
# keeps the cell formatting identical to the rest of the column.
$newRows = @(
    @{ A = "000061"; B = "Details: 60bpm_mf_44path_Normal_AllHabits. Script used: BasicVisualisation_AllAtOnce_HeatMap_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_SimpleCentroidTrackingData\Session01_Exp_F1_001_GHI_BlanksRemoved_SimpleCentroid.csv. File Location: Visualisations/Session01_SimpleCentroid_Figures. Date Generated: 25-Feb-2023 11:11:12" },
    @{ A = "000062"; B = "Details: 60bpm_mf_44path_Accelerando_AllHabits. Script used: BasicVisualisation_AllAtOnce_HeatMap_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_SimpleCentroidTrackingData\Session01_Exp_F2_001_GHI_BlanksRemoved_SimpleCentroid.csv. File Location: Visualisations/Session01_SimpleCentroid_Figures. Date Generated: 25-Feb-2023 11:11:22" },
    @{ A = "000063"; B = "Details: 60bpm_mf_44path_Ritardando_AllHabits. Script used: BasicVisualisation_AllAtOnce_HeatMap_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_SimpleCentroidTrackingData\Session01_Exp_F3_001_GHI_BlanksRemoved_SimpleCentroid.csv. File Location: Visualisations/Session01_SimpleCentroid_Figures. Date Generated: 25-Feb-2023 11:11:34" },
    @{ A = "000064"; B = "Details: 60bpm_mf_44path_Lead in_AllHabits. Script used: BasicVisualisation_AllAtOnce_HeatMap_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_SimpleCentroidTrackingData\Session01_Exp_F4_001_GHI_BlanksRemoved_SimpleCentroid.csv. File Location: Visualisations/Session01_SimpleCentroid_Figures. Date Generated: 25-Feb-2023 11:11:40" },
    @{ A = "000065"; B = "Details: 60bpm_mf_44path_Cut off_AllHabits. Script used: BasicVisualisation_AllAtOnce_HeatMap_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_SimpleCentroidTrackingData\Session01_Exp_F5_001_GHI_BlanksRemoved_SimpleCentroid.csv. File Location: Visualisations/Session01_SimpleCentroid_Figures. Date Generated: 25-Feb-2023 11:11:48" },
    @{ A = "000066"; B = "Details: 60bpm_mf_44path_Crescendo_AllHabits. Script used: BasicVisualisation_AllAtOnce_HeatMap_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_SimpleCentroidTrackingData\Session01_Exp_F6_001_GHI_BlanksRemoved_SimpleCentroid.csv. File Location: Visualisations/Session01_SimpleCentroid_Figures. Date Generated: 25-Feb-2023 11:12:00" }
)

$startRow = 62
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.A
    $cellA.Style = "Normal"

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $row.B
}
